# ---------------------------------------------------------------------------
# Applies the "ADDITIONAL SCRAPING" commit:
#   1. Transforms the MATCH_CARD_LINK URL columns on the existing
#      "ODI Batting" / "ODI Bowling" sheets into plain MATCH_CODE numbers,
#      renames the header, and drops the handful of truly-empty
#      INNING_NUMBER cells that used to be written out as empty strings.
#   2. Inserts a new "Player Info" sheet (before "ODI Batting") holding the
#      player's bio.
#   3. Appends a new "ODI Batting Extra" sheet (after "ODI Bowling") holding
#      additional per-match batting detail.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: rewrite a "...MatchCode=1234" URL cell to the bare "1234" code,
# keeping the cell as Text (not auto-coerced to a Number).
# ---------------------------------------------------------------------------
function Convert-MatchCodeColumn($ws, [int]$col, [int]$firstRow, [int]$lastRow) {
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $col)
        $val = $cell.Value()
        if ($val -ne $null -and "$val" -match 'MatchCode=(\d+)') {
            $code = $matches[1]
            $cell.NumberFormat = "@"
            $cell.Value = $code
        }
    }
}

# ---------------------------------------------------------------------------
# 1a. "ODI Batting" sheet
# ---------------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")

$lastRow = $batting.UsedRange.Rows.Count
$lastCol = $batting.UsedRange.Columns.Count

# D1: MATCH_CARD_LINK -> MATCH_CODE
$batting.Cells.Item(1, 4).Value = "MATCH_CODE"

# D2:Dn: full howstat URL -> bare match code, kept as text
Convert-MatchCodeColumn $batting 4 2 $lastRow

# Rows where INNING_NUMBER (col B) was only ever an empty placeholder string
# are now fully blank cells (no inlineStr at all).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $batting.Cells.Item($r, 2)
    $val = $cell.Value()
    if ($val -eq $null -or "$val" -eq "") {
        $cell.ClearContents()
    }
}

# ---------------------------------------------------------------------------
# 1b. "ODI Bowling" sheet
# ---------------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowlLastRow = $bowling.UsedRange.Rows.Count

# B1: MATCH_CARD_LINK -> MATCH_CODE
$bowling.Cells.Item(1, 2).Value = "MATCH_CODE"

# B2:Bn: full howstat URL -> bare match code, kept as text
Convert-MatchCodeColumn $bowling 2 2 $bowlLastRow

# ---------------------------------------------------------------------------
# 2. New "Player Info" sheet, inserted before "ODI Batting"
# ---------------------------------------------------------------------------
$batting.Activate()
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $piHeaders.Length; $c++) {
    $cell = $playerInfo.Cells.Item(1, $c)
    $cell.Value = $piHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$piRow = @("3348", "Luteru Ross Poutoa Lote Taylor", "Right Handed", "Right Arm Off Break")
for ($c = 1; $c -le $piRow.Length; $c++) {
    $cell = $playerInfo.Cells.Item(2, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $piRow[$c - 1]
}

# ---------------------------------------------------------------------------
# 3. New "ODI Batting Extra" sheet, appended after "ODI Bowling"
# ---------------------------------------------------------------------------
# Re-resolve the "ODI Bowling" sheet now (its position shifted after the
# "Player Info" insert above), otherwise Add() anchors off a stale sheet.
$bowlingNow = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Add($null, $bowlingNow)
$extra.Name = "ODI Batting Extra"

$exHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $exHeaders.Length; $c++) {
    $cell = $extra.Cells.Item(1, $c)
    $cell.Value = $exHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$exData = @(
    @("4251","4","3","0","9.17%","NO"),
    @("4252","4","7","0","20.91%","NO"),
    @("4305","4","","","","NO"),
    @("4311","","","","","NO"),
    @("4315","4","6","1","27.75%","NO"),
    @("4328","4","0","0","0.41%","NO"),
    @("4333","","","","","NO"),
    @("4337","","","","","NO"),
    @("4341","4","2","0","19.11%","NO"),
    @("4346","","","","","NO"),
    @("4353","","","","","NO"),
    @("4355","","","","","NO"),
    @("4402","","","","","NO"),
    @("4406","","","","","NO"),
    @("4410","4","1","0","4.00%","NO"),
    @("4423","4","0","0","2.14%","NO"),
    @("4455","4","1","0","2.20%","NO"),
    @("4563","4","1","0","5.39%","NO"),
    @("4566","4","0","0","0.38%","NO"),
    @("4568","","","","","NO")
)

for ($i = 0; $i -lt $exData.Length; $i++) {
    $r = $i + 2
    $row = $exData[$i]

    # A: MATCH_CODE - text
    $cellA = $extra.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[0]

    # B: BATTING_POSITION - numeric when present; an explicitly-blank cell
    # (not simply absent) when not, matching the source export's behaviour.
    $cellB = $extra.Cells.Item($r, 2)
    $cellB.NumberFormat = "General"
    if ($row[1] -ne "") {
        $cellB.Value = [double]$row[1]
    } else {
        $cellB.Value = ""
    }

    # C: NUM_4 - text
    $cellC = $extra.Cells.Item($r, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $row[2]

    # D: NUM_6 - text
    $cellD = $extra.Cells.Item($r, 4)
    $cellD.NumberFormat = "@"
    $cellD.Value = $row[3]

    # E: PERCENT_RUNS_OF_TOTAL - text
    $cellE = $extra.Cells.Item($r, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $row[4]

    # F: MAN_OF_MATCH - text
    $cellF = $extra.Cells.Item($r, 6)
    $cellF.NumberFormat = "@"
    $cellF.Value = $row[5]
}


